$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 57, shifting rows 57:64 down to 58:65
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with its values
$ws.Range("A57").Value = 7
$ws.Range("B57").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C57").Value = "Ñuble"
$ws.Range("D57").Value = 44610
$ws.Range("E57").Value = 16
$ws.Range("F57").Value = 100112022
$ws.Range("G57").Value = "Arveja Verde"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 60
$ws.Range("K57").Value = 24000
$ws.Range("L57").Value = 25000
$ws.Range("M57").Value = 24500
$ws.Range("N57").Value = "$/saco 25 kilos"
$ws.Range("O57").Value = "Provincia de Diguillín"
$ws.Range("P57").Value = 980
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"
